# "Generate Report for Handback"
#
# The handback for a5506523-...-zh-cn / de-de finished, so the status that
# used to read "Ready for handoff" is now "Handed back: in sync with en-US"
# everywhere it is shown (the Overview roll-up columns as well as the
# per-language detail sheets), the "Latest Handback DateTime" for the
# first (a5506523...) row on each language sheet is refreshed to the
# handback time, and the stale "handback file is not the latest" error
# message is cleared out now that the handback is in sync.
#
# Columns that show this text auto-widen/narrow to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status roll-up columns (E, F) for both
# file rows (2, 3).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Latest Handback DateTime for the a5506523... file just got generated.
$zhcn.Range("K2").Value = "2016-10-18 05:20:07"
$zhcn.Range("K3").Value = "2016-10-18 05:20:07"

# The handback is now in sync, so the stale "not the latest" error detail
# for the a5506523... row is cleared.
$zhcn.Range("P2").Value = ""

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("K2").Value = "2016-10-18 05:20:31"
$dede.Range("K3").Value = "2016-10-18 05:20:31"

$dede.Range("P2").Value = ""

# ---------------------------------------------------------------------
# Column widths follow the new (longer) status text / (now empty) error
# detail text.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 12.75

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(16).ColumnWidth = 12.75
